$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to their rounded (2-decimal) equivalents
$ws.Range("B5").Value = 0.91
$ws.Range("C5").Value = 0.43
$ws.Range("D5").Value = 0.58
$ws.Range("E5").Value = 2.22
$ws.Range("F5").Value = 1.1
$ws.Range("G5").Value = 0.66
$ws.Range("H5").Value = 9.300000000000001
$ws.Range("I5").Value = 1.16
$ws.Range("J5").Value = 0.53
$ws.Range("K5").Value = 0.34
$ws.Range("L5").Value = 0.82
$ws.Range("M5").Value = 0.92
$ws.Range("N5").Value = 0.28
$ws.Range("O5").Value = 0.75
$ws.Range("P5").Value = 1.15
$ws.Range("Q5").Value = 0.91
$ws.Range("R5").Value = 0.65
$ws.Range("S5").Value = 0.26
$ws.Range("T5").Value = 4.28
$ws.Range("U5").Value = 2.64
$ws.Range("V5").Value = 0.6899999999999999
$ws.Range("W5").Value = 1.71
$ws.Range("X5").Value = 0.75
$ws.Range("Y5").Value = 0.42
$ws.Range("Z5").Value = 4.11
$ws.Range("AA5").Value = 0.61
$ws.Range("AB5").Value = 0.7
$ws.Range("AC5").Value = 0.79
$ws.Range("AD5").Value = 0.7
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 9.16
$ws.Range("AG5").Value = 0.25
$ws.Range("AH5").Value = 0.88

# Delete row 6 entirely
$ws.Rows.Item(6).Delete()
